# Insert a new data row at row 296 (pushing existing rows 296-349 down to 297-350)
# and populate it with a new "Pepino ensalada" observation for Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(296).Insert()

$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = 44694
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 100112043
$ws.Range("G296").Value = "Pepino ensalada"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 115
$ws.Range("K296").Value = 17000
$ws.Range("L296").Value = 18000
$ws.Range("M296").Value = 17478
$ws.Range("N296").Value = "$/caja 70 unidades"
$ws.Range("O296").Value = "Región de Arica y Parinacota"
$ws.Range("P296").Value = 250
$ws.Range("Q296").Value = 70
$ws.Range("R296").Value = "Hortaliza"
